$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entidade")

# Row 2 ("entidade" record) field edits
$ws.Range("C2").Value = "07.046.881/1007-32"
$ws.Range("E2").Value = 141

# Leave the cursor where the author left it when they saved
$ws.Range("E6:F6").Select()
